$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 19477
$ws.Range("B3").Value = 14574
$ws.Range("B4").Value = 1818
$ws.Range("B5").Value = 18042
